# The deck's theme (ppt/theme/theme1.xml, used by the one SlideMaster /
# all layouts / all slides) was switched from the custom "Integral" theme
# palette to the built-in "Office Theme" palette (dk1/lt1/dk2/lt2/accent1-6/
# hlink/folHlink). Font scheme and format scheme are identical between the
# two themes, so only the 12 theme colors need to change.
#
# We drive this the same way a user would from the Design tab - by setting
# the presentation's theme colors - via the ThemeColorScheme object (the
# modern 12-slot equivalent of the legacy 8-slot Master.ColorScheme), which
# PowerPoint's object model exposes per slide / slide range and which maps
# back onto the shared slide-master theme part.

$p = $ppt.ActivePresentation

function Set-ThemeColor {
    param($Scheme, [int]$Index, [string]$Hex)
    $r = [Convert]::ToInt32($Hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4,2), 16)
    # PowerPoint's RGB property packs colour as R + G*256 + B*65536
    $Scheme.Colors($Index).RGB = $r + ($g * 256) + ($b * 65536)
}

# All slides share the single slide master, so changing the theme colours
# through any one slide's ThemeColorScheme updates the shared theme part
# (ppt/theme/theme1.xml) used by the whole deck.
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

Set-ThemeColor $themeColors 1  "000000"   # dk1
Set-ThemeColor $themeColors 2  "FFFFFF"   # lt1
Set-ThemeColor $themeColors 3  "44546A"   # dk2
Set-ThemeColor $themeColors 4  "E7E6E6"   # lt2
Set-ThemeColor $themeColors 5  "5B9BD5"   # accent1
Set-ThemeColor $themeColors 6  "ED7D31"   # accent2
Set-ThemeColor $themeColors 7  "A5A5A5"   # accent3
Set-ThemeColor $themeColors 8  "FFC000"   # accent4
Set-ThemeColor $themeColors 9  "4472C4"   # accent5
Set-ThemeColor $themeColors 10 "70AD47"   # accent6
Set-ThemeColor $themeColors 11 "0563C1"   # hlink
Set-ThemeColor $themeColors 12 "954F72"   # folHlink
